$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 74 (this shifts the existing "347-Top K Frequent
# Elements" / "271 - Encode and Decode Strings" rows, and the trailing
# blank styled rows, down by one -- matching the table growing from
# A2:X84 to A2:X85).
$ws.Rows("74:74").Insert()

# Populate the newly solved LeetCode problem: 238 - Product of Array
# Except Self.
$ws.Range("A74").Value = "Array"
$ws.Range("B74").Value = 238
$ws.Range("C74").Value = "238-Product of Array Except Self"
$ws.Range("D74").Value = "Medium"
$ws.Range("E74").Value = "Store left product and right product arrays to allow solution to be computed in linear time"
$ws.Range("F74").Value = "O(n) time, O(n) memory"
$ws.Range("G74").Value = "O(n) time, O(1) memory"
$ws.Range("H74").Value = "Store left product and right product arrays within the input and output arrays"
$ws.Range("I74").Value = "O(1) memory"
$ws.Range("J74").Value = "no"
$ws.Range("K74").Value = "yes"
$ws.Range("M74").Value = "2 hours 30 minutes"
$ws.Range("N74").Value = "yes"

# Matches the wrap-text style used elsewhere for the longer free-text
# "algorithm description" columns (E and H).
$ws.Range("E74").WrapText = $true
$ws.Range("H74").WrapText = $true

$ws.Rows("74:74").RowHeight = 29

# The row insert grew the used range (A2:X84 -> A2:X85) but the table
# definition itself needs to be told to cover the extra row explicitly.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:X85"))

# Update the view to reflect where the author was working / had selected.
$excel.ActiveWindow.ScrollRow = 70
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I73").Select()
